$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(131).Insert()
$ws.Rows.Item(132).Copy()
$ws.Rows.Item(131).PasteSpecial(-4122)
Write-Host "done"
